$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.777961
$ws.Range("H2").Value = 2.333883
$ws.Range("I2").Value = 0.7646397019917995
$ws.Range("J2").Value = 0.7646397019917995
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5535766666666666
$ws.Range("N2").Value = 1.66073
$ws.Range("O2").Value = 0.7546693956119359
$ws.Range("P2").Value = 0.754669395611936
$ws.Range("Q2").Value = 0.4306610571766666
$ws.Range("R2").Value = 3.87594951459
$ws.Range("S2").Value = 0.5770501817630421
$ws.Range("T2").Value = 0.5770501817630422

# Row 3
$ws.Range("G3").Value = 0.777961
$ws.Range("H3").Value = 2.333883
$ws.Range("I3").Value = 0.7646397019917995
$ws.Range("J3").Value = 0.7646397019917995
$ws.Range("M3").Value = 0.1799586666666667
$ws.Range("N3").Value = 0.539876
$ws.Range("O3").Value = 0.245330604388064
$ws.Range("P3").Value = 0.2453306043880641
$ws.Range("Q3").Value = 0.1400008242786667
$ws.Range("R3").Value = 1.260007418508
$ws.Range("S3").Value = 0.1875895202287574
$ws.Range("T3").Value = 0.1875895202287574

# Row 4
$ws.Range("I4").Value = 0.2353602980082005
$ws.Range("J4").Value = 0.2353602980082005
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5535766666666666
$ws.Range("N4").Value = 1.66073
$ws.Range("O4").Value = 0.7546693956119359
$ws.Range("P4").Value = 0.754669395611936
$ws.Range("Q4").Value = 0.1325598376511111
$ws.Range("R4").Value = 1.19303853886
$ws.Range("S4").Value = 0.1776192138488938
$ws.Range("T4").Value = 0.1776192138488938

# Row 5
$ws.Range("I5").Value = 0.2353602980082005
$ws.Range("J5").Value = 0.2353602980082005
$ws.Range("M5").Value = 0.1799586666666667
$ws.Range("N5").Value = 0.539876
$ws.Range("O5").Value = 0.245330604388064
$ws.Range("P5").Value = 0.2453306043880641
$ws.Range("Q5").Value = 0.04309302229244445
$ws.Range("R5").Value = 0.3878372006320001
$ws.Range("S5").Value = 0.05774108415930669
$ws.Range("T5").Value = 0.05774108415930669
